$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "唉, 找回来了"
$ws.Range("D4").Select()
